# The sheet originally has 10 data rows (rows 2-11) below the header (row 1).
# Two blank separator rows are inserted into the data table:
#   - one before the old row 7 (pushing old rows 7-11 down to 8-12)
#   - one before what is now row 11 (pushing the former old rows 10-11,
#     now at 11-12, down to 12-13)
# The net effect matches the target: rows 1-6 stay put, a blank row appears
# at row 7, the next three records shift to rows 8-10, another blank row
# appears at row 11, and the final two records land on rows 12-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Insert() | Out-Null
$ws.Rows.Item(11).Insert() | Out-Null

$ws.Range("C20").Select() | Out-Null
